$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "/home/daniel/Spike Data/Matlab files/Exp 19 baseline data.mat"
$ws.Range("A3").Value = "/home/daniel/Spike Data/Matlab files/Exp 27 unit 1 data.mat"
$ws.Range("A4").Value = "/home/daniel/Spike Data/Matlab files/Exp 27 unit 2 data.mat"

$ws.Range("F26").Select()
